#
# Fruta / hortaliza, semanal
#
# Inserts a new weekly data row at row 3 (pushing the existing rows 3-13
# down to 4-14), and populates it with the new week's readings for
# "Terminal La Palmera de La Serena - Arándano (blue)".
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 3; this shifts rows 3..13 down
# to 4..14 and grows the used range to A1:T14 (matches Excel's native
# "format from above" behaviour, carrying the date-format style on column D).
$ws.Rows("3").Insert()

# Populate the newly inserted row 3 with this week's record.
$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "Terminal La Palmera de La Serena"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44490
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100101
$ws.Range("H3").Value = "Berries"
$ws.Range("I3").Value = 100101001
$ws.Range("J3").Value = "Arándano (blue)"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 160
$ws.Range("N3").Value = 11500
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 11750
$ws.Range("Q3").Value = "$/bandeja 2 kilos"
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 5875
$ws.Range("T3").Value = 2
